# The workbook has two sheets: "Users list - M2M" (active/selected tab)
# and "Companies (Ignore one)". The edit only touches the first sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 used to hold a (bogus) numeric id of 2 in column A, matching no
# real record. Replace it with a text id that clearly doesn't exist so the
# import test fixture fails for the expected reason ("ThatIdDoesntExist").
$ws.Range("A3").Value = "ThatIdDoesntExist"

# Move the visible selection/active cell to A2.
$ws.Range("A2").Select()
